# Auto-generated edit script applying the value changes described by the diff.
# Each worksheet's cells are updated to their new literal values; cells removed
# in the diff are cleared (which drops the <c> element on save).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 545.2222
$ws.Range("I11").Value = 545.2222
$ws.Range("K11").Value = 545.2222
$ws.Range("M11").Value = -405.2222
$ws.Range("H18").Value = 986.75
$ws.Range("I18").Value = 999
$ws.Range("J18").Value = 950
$ws.Range("K18").Value = 999
$ws.Range("L18").Value = 950
$ws.Range("M18").Value = -715
$ws.Range("N18").Value = -1518
$ws.Range("H40").Value = 2809
$ws.Range("I40").Value = 2306.6
$ws.Range("J40").Value = 3227.6667
$ws.Range("K40").Value = 2306.6
$ws.Range("L40").Value = 3227.6667
$ws.Range("M40").Value = -2131.6
$ws.Range("N40").Value = -3577.6667
$ws.Range("H113").Value = 5424.5356
$ws.Range("J113").Value = 7368.5
$ws.Range("L113").Value = 7368.5
$ws.Range("N113").Value = -13876.5
$ws.Range("H116").Value = 6056.5557
$ws.Range("I116").Value = 4930
$ws.Range("K116").Value = 4930
$ws.Range("M116").Value = -1488
$ws.Range("H137").Value = 2938.8948
$ws.Range("I137").Value = 2226.111
$ws.Range("K137").Value = 6678.333
$ws.Range("M137").Value = -4128.333
$ws.Range("H138").Value = 2371.23
$ws.Range("I138").Value = 1880.4412
$ws.Range("J138").Value = 2624.0605
$ws.Range("K138").Value = 5641.3236
$ws.Range("L138").Value = 7872.181500000001
$ws.Range("M138").Value = -501.3235999999997
$ws.Range("N138").Value = -18152.1815
$ws.Range("H141").Value = 3341.6052
$ws.Range("I141").Value = 2211.697
$ws.Range("J141").Value = 10799
$ws.Range("K141").Value = 6635.091
$ws.Range("L141").Value = 32397
$ws.Range("M141").Value = -1455.091
$ws.Range("N141").Value = -42757

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1253.2609
$ws.Range("I2").Value = 697.8333
$ws.Range("J2").Value = 3252.8
$ws.Range("K2").Value = 697.8333
$ws.Range("L2").Value = 3252.8
$ws.Range("M2").Value = -584.8333
$ws.Range("N2").Value = -3478.8
$ws.Range("H61").Value = 2884.56
$ws.Range("I61").Value = 2884.56
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2884.56
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = $null
$ws.Range("H110").Value = 1461.875
$ws.Range("I110").Value = 1465.421
$ws.Range("J110").Value = 1448.4
$ws.Range("K110").Value = 1465.421
$ws.Range("L110").Value = 1448.4
$ws.Range("M110").Value = 579.579
$ws.Range("N110").Value = -5538.4
$ws.Range("H116").Value = 1253.2609
$ws.Range("I116").Value = 697.8333
$ws.Range("J116").Value = 3252.8
$ws.Range("K116").Value = 697.8333
$ws.Range("L116").Value = 3252.8
$ws.Range("M116").Value = 1596.1667
$ws.Range("N116").Value = -7840.8
$ws.Range("H132").Value = 1327
$ws.Range("I132").Value = 1264.8372
$ws.Range("K132").Value = 3794.5116
$ws.Range("M132").Value = -1264.5116
$ws.Range("H136").Value = 2884.56
$ws.Range("I136").Value = 2884.56
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8653.68
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1253.2609
$ws.Range("I3").Value = 697.8333
$ws.Range("J3").Value = 3252.8
$ws.Range("K3").Value = 697.8333
$ws.Range("L3").Value = 3252.8
$ws.Range("M3").Value = -583.8333
$ws.Range("N3").Value = -3480.8
$ws.Range("H86").Value = 3620.625
$ws.Range("I86").Value = 2657.3333
$ws.Range("J86").Value = 4198.6
$ws.Range("K86").Value = 2657.3333
$ws.Range("L86").Value = 4198.6
$ws.Range("M86").Value = -1534.3333
$ws.Range("N86").Value = -6444.6
$ws.Range("H89").Value = 3620.625
$ws.Range("I89").Value = 2657.3333
$ws.Range("J89").Value = 4198.6
$ws.Range("K89").Value = 13286.6665
$ws.Range("L89").Value = 20993
$ws.Range("M89").Value = -7670.666499999999
$ws.Range("N89").Value = -32225
$ws.Range("H134").Value = 813.1
$ws.Range("I134").Value = 599.75
$ws.Range("K134").Value = 1799.25
$ws.Range("M134").Value = 735.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 10999.944
$ws.Range("I41").Value = 9000
$ws.Range("J41").Value = 11249.9375
$ws.Range("K41").Value = 9000
$ws.Range("L41").Value = 11249.9375
$ws.Range("M41").Value = -8572
$ws.Range("N41").Value = -12105.9375
$ws.Range("H59").Value = 46359
$ws.Range("J59").Value = 47994.555
$ws.Range("L59").Value = 47994.555
$ws.Range("N59").Value = -50284.555
$ws.Range("H60").Value = 14166.667
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").Value = $null
$ws.Range("H62").Value = 8199
$ws.Range("I62").Value = 8124.5
$ws.Range("K62").Value = 8124.5
$ws.Range("M62").Value = -7500.5
$ws.Range("H65").Value = 8199
$ws.Range("I65").Value = 8124.5
$ws.Range("K65").Value = 40622.5
$ws.Range("M65").Value = -37502.5
$ws.Range("H68").Value = 62854.89
$ws.Range("J68").Value = 62854.89
$ws.Range("L68").Value = 62854.89
$ws.Range("N68").Value = -64352.89
$ws.Range("H69").Value = 33960.668
$ws.Range("I69").Value = 33441
$ws.Range("J69").Value = 35000
$ws.Range("K69").Value = 33441
$ws.Range("L69").Value = 35000
$ws.Range("M69").Value = -32692
$ws.Range("N69").Value = -36498
$ws.Range("H71").Value = 62854.89
$ws.Range("J71").Value = 62854.89
$ws.Range("L71").Value = 188564.67
$ws.Range("N71").Value = -196052.67
$ws.Range("H72").Value = 33960.668
$ws.Range("I72").Value = 33441
$ws.Range("J72").Value = 35000
$ws.Range("K72").Value = 100323
$ws.Range("L72").Value = 105000
$ws.Range("M72").Value = -96579
$ws.Range("N72").Value = -112488
$ws.Range("H132").Value = 848.44446
$ws.Range("I132").Value = 591
$ws.Range("J132").Value = 1749.5
$ws.Range("K132").Value = 1773
$ws.Range("L132").Value = 5248.5
$ws.Range("M132").Value = 757
$ws.Range("N132").Value = -10308.5
$ws.Range("H134").Value = 1770.8379
$ws.Range("I134").Value = 1815.8182
$ws.Range("K134").Value = 5447.4546
$ws.Range("M134").Value = -2912.4546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 217.21053
$ws.Range("I12").Value = 304.5
$ws.Range("J12").Value = 153.72728
$ws.Range("K12").Value = 913.5
$ws.Range("L12").Value = 461.18184
$ws.Range("M12").Value = -740.5
$ws.Range("N12").Value = -807.18184
$ws.Range("H113").Value = 552.7646999999999
$ws.Range("I113").Value = 510.85715
$ws.Range("J113").Value = 582.1
$ws.Range("K113").Value = 1532.57145
$ws.Range("L113").Value = 1746.3
$ws.Range("M113").Value = 637.4285500000001
$ws.Range("N113").Value = -6086.3
$ws.Range("H114").Value = 2579.6
$ws.Range("J114").Value = 2481.1667
$ws.Range("L114").Value = 7443.500100000001
$ws.Range("N114").Value = -13951.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 19846.8
$ws.Range("J95").Value = 19846.8
$ws.Range("L95").Value = 19846.8
$ws.Range("N95").Value = -25338.8
$ws.Range("H132").Value = 1360.7587
$ws.Range("I132").Value = 790.4583
$ws.Range("K132").Value = 2371.3749
$ws.Range("M132").Value = 158.6251000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2824.3333
$ws.Range("I7").Value = 1625
$ws.Range("K7").Value = 1625
$ws.Range("M7").Value = -1513
$ws.Range("H46").Value = 2407.3572
$ws.Range("J46").Value = 3063.125
$ws.Range("L46").Value = 3063.125
$ws.Range("N46").Value = -3439.125
$ws.Range("H68").Value = 3045.3333
$ws.Range("I68").Value = 3065.6667
$ws.Range("K68").Value = 3065.6667
$ws.Range("M68").Value = -2316.6667
$ws.Range("H71").Value = 3045.3333
$ws.Range("I71").Value = 3065.6667
$ws.Range("K71").Value = 15328.3335
$ws.Range("M71").Value = -11584.3335
$ws.Range("H82").Value = 495.3
$ws.Range("H85").Value = 495.3
$ws.Range("H126").Value = 2824.3333
$ws.Range("I126").Value = 1625
$ws.Range("K126").Value = 4875
$ws.Range("M126").Value = -2405
$ws.Range("H136").Value = 2268.7307
$ws.Range("I136").Value = 2308.0908
$ws.Range("K136").Value = 6924.2724
$ws.Range("M136").Value = -4374.2724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 1440
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = $null
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").Value = $null
$ws.Range("H14").Value = 4990
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").Value = $null
$ws.Range("H113").Value = 1044.8889
$ws.Range("I113").Value = 1121.4286
$ws.Range("K113").Value = 3364.2858
$ws.Range("M113").Value = -1194.2858
$ws.Range("H126").Value = 500.5
$ws.Range("I126").Value = 500.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 1501.5
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = $null

